$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet1 "Recommandations": update C/D/E for rows 2-23 (titles/F/G unchanged) ---
$ws1.Range("C2").Value = 5
$ws1.Range("D2").Value = 4850
$ws1.Range("E2").Value = 980
$ws1.Range("C3").Value = 10
$ws1.Range("D3").Value = 4234.51
$ws1.Range("E3").Value = 112.08
$ws1.Range("C4").Value = 5
$ws1.Range("D4").Value = 3430
$ws1.Range("E4").Value = 695
$ws1.Range("C5").Value = 5
$ws1.Range("D5").Value = 3400
$ws1.Range("E5").Value = 685
$ws1.Range("C6").Value = 5
$ws1.Range("D6").Value = 3307.48
$ws1.Range("E6").Value = 666.3200000000001
$ws1.Range("C7").Value = 5
$ws1.Range("D7").Value = 2965
$ws1.Range("E7").Value = 600
$ws1.Range("C8").Value = 5
$ws1.Range("D8").Value = 2925
$ws1.Range("E8").Value = 585
$ws1.Range("C9").Value = 5
$ws1.Range("D9").Value = 2825
$ws1.Range("E9").Value = 580
$ws1.Range("C10").Value = 5
$ws1.Range("D10").Value = 2690
$ws1.Range("E10").Value = 540
$ws1.Range("C11").Value = 5
$ws1.Range("D11").Value = 1872.35
$ws1.Range("E11").Value = 374.85
$ws1.Range("C12").Value = 5
$ws1.Range("D12").Value = 1745.24
$ws1.Range("E12").Value = 350.03
$ws1.Range("C13").Value = 5
$ws1.Range("D13").Value = 1664.04
$ws1.Range("E13").Value = 333.07
$ws1.Range("C14").Value = 5
$ws1.Range("D14").Value = 986.76
$ws1.Range("E14").Value = 203.97
$ws1.Range("C15").Value = 5
$ws1.Range("D15").Value = 882.98
$ws1.Range("E15").Value = 177.85
$ws1.Range("C16").Value = 5
$ws1.Range("D16").Value = 874.46
$ws1.Range("E16").Value = 178.9
$ws1.Range("C17").Value = 5
$ws1.Range("D17").Value = 657.97
$ws1.Range("E17").Value = 132.12
$ws1.Range("C18").Value = 5
$ws1.Range("D18").Value = 645.72
$ws1.Range("E18").Value = 133.24
$ws1.Range("C19").Value = 5
$ws1.Range("D19").Value = 612.92
$ws1.Range("E19").Value = 122.75
$ws1.Range("C20").Value = 5
$ws1.Range("D20").Value = 602.37
$ws1.Range("E20").Value = 120.64
$ws1.Range("C21").Value = 5
$ws1.Range("D21").Value = 555.77
$ws1.Range("E21").Value = 111.84
$ws1.Range("C22").Value = 5
$ws1.Range("D22").Value = 532.39
$ws1.Range("E22").Value = 106.19
$ws1.Range("C23").Value = 5
$ws1.Range("D23").Value = 477.87
$ws1.Range("E23").Value = 95.34

# --- Sheet1 "Recommandations": rewrite full rows 24-54 (reordered / new companies) ---
$ws1.Range("A24").Value = 'SOLIBRA CI (SLBC)'
$ws1.Range("B24").Value = 2
$ws1.Range("C24").Value = 0
$ws1.Range("D24").Value = 11.44
$ws1.Range("E24").Value = 7.48
$ws1.Range("F24").Value = '🟡 Observer'
$ws1.Range("G24").Value = '➖ Neutre'
$ws1.Range("A25").Value = 'BERNABE CI (BNBC)'
$ws1.Range("B25").Value = 3
$ws1.Range("C25").Value = 1
$ws1.Range("D25").Value = 10
$ws1.Range("E25").Value = 7.32
$ws1.Range("F25").Value = '🟢 Achat'
$ws1.Range("G25").Value = '✅ Renforcer'
$ws1.Range("A26").Value = 'FILTISAC CI (FTSC)'
$ws1.Range("B26").Value = 2
$ws1.Range("C26").Value = 1
$ws1.Range("D26").Value = 7.57
$ws1.Range("E26").Value = 7.47
$ws1.Range("F26").Value = '🟡 Observer'
$ws1.Range("G26").Value = '👀 À surveiller'
$ws1.Range("A27").Value = 'UNILEVER CI (UNLC)'
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 0
$ws1.Range("D27").Value = 7.49
$ws1.Range("E27").Value = 7.49
$ws1.Range("F27").Value = '🟡 Observer'
$ws1.Range("G27").Value = '➖ Neutre'
$ws1.Range("A28").Value = 'SONATEL SN (SNTS)'
$ws1.Range("B28").Value = 2
$ws1.Range("C28").Value = 0
$ws1.Range("D28").Value = 7.2
$ws1.Range("E28").Value = 3.59
$ws1.Range("F28").Value = '🟡 Observer'
$ws1.Range("G28").Value = '➖ Neutre'
$ws1.Range("A29").Value = 'CORIS BANK INTERNATIONAL (CBIBF)'
$ws1.Range("B29").Value = 1
$ws1.Range("C29").Value = 0
$ws1.Range("D29").Value = 7.18
$ws1.Range("E29").Value = 7.18
$ws1.Range("F29").Value = '🟡 Observer'
$ws1.Range("G29").Value = '➖ Neutre'
$ws1.Range("A30").Value = 'SICOR CI (SICC)'
$ws1.Range("B30").Value = 1
$ws1.Range("C30").Value = 0
$ws1.Range("D30").Value = 7.1
$ws1.Range("E30").Value = 7.1
$ws1.Range("F30").Value = '🟡 Observer'
$ws1.Range("G30").Value = '➖ Neutre'
$ws1.Range("A31").Value = 'BANK OF AFRICA ML (BOAM)'
$ws1.Range("B31").Value = 1
$ws1.Range("C31").Value = 0
$ws1.Range("D31").Value = 6.22
$ws1.Range("E31").Value = 6.22
$ws1.Range("F31").Value = '🟡 Observer'
$ws1.Range("G31").Value = '➖ Neutre'
$ws1.Range("A32").Value = 'BANK OF AFRICA SENEGAL (BOAS)'
$ws1.Range("B32").Value = 1
$ws1.Range("C32").Value = 0
$ws1.Range("D32").Value = 3.38
$ws1.Range("E32").Value = 3.38
$ws1.Range("F32").Value = '🟡 Observer'
$ws1.Range("G32").Value = '➖ Neutre'
$ws1.Range("A33").Value = 'SUCRIVOIRE (SCRC)'
$ws1.Range("B33").Value = 1
$ws1.Range("C33").Value = 0
$ws1.Range("D33").Value = 3.16
$ws1.Range("E33").Value = 3.16
$ws1.Range("F33").Value = '🟡 Observer'
$ws1.Range("G33").Value = '➖ Neutre'
$ws1.Range("A34").Value = 'SERVAIR ABIDJAN CI (ABJC)'
$ws1.Range("B34").Value = 1
$ws1.Range("C34").Value = 0
$ws1.Range("D34").Value = 2.99
$ws1.Range("E34").Value = 2.99
$ws1.Range("F34").Value = '🟡 Observer'
$ws1.Range("G34").Value = '➖ Neutre'
$ws1.Range("A35").Value = 'PALM CI (PALC)'
$ws1.Range("B35").Value = 1
$ws1.Range("C35").Value = 1
$ws1.Range("D35").Value = 2.01
$ws1.Range("E35").Value = -2.74
$ws1.Range("F35").Value = '🟡 Observer'
$ws1.Range("G35").Value = '👀 À surveiller'
$ws1.Range("A36").Value = 'SODE CI (SDCC)'
$ws1.Range("B36").Value = 1
$ws1.Range("C36").Value = 0
$ws1.Range("D36").Value = 1.67
$ws1.Range("E36").Value = 1.67
$ws1.Range("F36").Value = '🟡 Observer'
$ws1.Range("G36").Value = '➖ Neutre'
$ws1.Range("A37").Value = 'ONATEL BF (ONTBF)'
$ws1.Range("B37").Value = 1
$ws1.Range("C37").Value = 0
$ws1.Range("D37").Value = 1.35
$ws1.Range("E37").Value = 1.35
$ws1.Range("F37").Value = '🟡 Observer'
$ws1.Range("G37").Value = '➖ Neutre'
$ws1.Range("A38").Value = 'VIVO ENERGY CI (SHEC)'
$ws1.Range("B38").Value = 1
$ws1.Range("C38").Value = 0
$ws1.Range("D38").Value = 0.9399999999999999
$ws1.Range("E38").Value = 0.9399999999999999
$ws1.Range("F38").Value = '🟡 Observer'
$ws1.Range("G38").Value = '➖ Neutre'
$ws1.Range("A39").Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws1.Range("B39").Value = 1
$ws1.Range("C39").Value = 1
$ws1.Range("D39").Value = 0.39
$ws1.Range("E39").Value = 3.26
$ws1.Range("F39").Value = '🟡 Observer'
$ws1.Range("G39").Value = '👀 À surveiller'
$ws1.Range("A40").Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws1.Range("B40").Value = 1
$ws1.Range("C40").Value = 1
$ws1.Range("D40").Value = 0.32
$ws1.Range("E40").Value = -5.56
$ws1.Range("F40").Value = '🟡 Observer'
$ws1.Range("G40").Value = '👀 À surveiller'
$ws1.Range("A41").Value = 'TOTAL'
$ws1.Range("B41").Value = 0
$ws1.Range("C41").Value = 5
$ws1.Range("D41").Value = 0
$ws1.Range("E41").Value = 0
$ws1.Range("F41").Value = '🟡 Observer'
$ws1.Range("G41").Value = '➖ Neutre'
$ws1.Range("A42").Value = 'SOGB CI (SOGC)'
$ws1.Range("B42").Value = 0
$ws1.Range("C42").Value = 1
$ws1.Range("D42").Value = -1.22
$ws1.Range("E42").Value = -1.22
$ws1.Range("F42").Value = '🟡 Observer'
$ws1.Range("G42").Value = '➖ Neutre'
$ws1.Range("A43").Value = 'SAFCA CI (SAFC)'
$ws1.Range("B43").Value = 2
$ws1.Range("C43").Value = 1
$ws1.Range("D43").Value = -1.25
$ws1.Range("E43").Value = 3.08
$ws1.Range("F43").Value = '🟡 Observer'
$ws1.Range("G43").Value = '👀 À surveiller'
$ws1.Range("A44").Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$ws1.Range("B44").Value = 0
$ws1.Range("C44").Value = 1
$ws1.Range("D44").Value = -1.4
$ws1.Range("E44").Value = -1.4
$ws1.Range("F44").Value = '🟡 Observer'
$ws1.Range("G44").Value = '➖ Neutre'
$ws1.Range("A45").Value = 'BICI CI (BICC)'
$ws1.Range("B45").Value = 0
$ws1.Range("C45").Value = 1
$ws1.Range("D45").Value = -2.6
$ws1.Range("E45").Value = -2.6
$ws1.Range("F45").Value = '🟡 Observer'
$ws1.Range("G45").Value = '➖ Neutre'
$ws1.Range("A46").Value = 'BANK OF AFRICA CI (BOAC)'
$ws1.Range("B46").Value = 0
$ws1.Range("C46").Value = 1
$ws1.Range("D46").Value = -3.27
$ws1.Range("E46").Value = -3.27
$ws1.Range("F46").Value = '🟡 Observer'
$ws1.Range("G46").Value = '➖ Neutre'
$ws1.Range("A47").Value = 'ORANGE COTE D''IVOIRE (ORAC)'
$ws1.Range("B47").Value = 0
$ws1.Range("C47").Value = 1
$ws1.Range("D47").Value = -3.33
$ws1.Range("E47").Value = -3.33
$ws1.Range("F47").Value = '🟡 Observer'
$ws1.Range("G47").Value = '➖ Neutre'
$ws1.Range("A48").Value = 'NEI-CEDA CI (NEIC)'
$ws1.Range("B48").Value = 0
$ws1.Range("C48").Value = 1
$ws1.Range("D48").Value = -3.33
$ws1.Range("E48").Value = -3.33
$ws1.Range("F48").Value = '🟡 Observer'
$ws1.Range("G48").Value = '➖ Neutre'
$ws1.Range("A49").Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$ws1.Range("B49").Value = 1
$ws1.Range("C49").Value = 2
$ws1.Range("D49").Value = -3.7
$ws1.Range("E49").Value = -3.85
$ws1.Range("F49").Value = '🟡 Observer'
$ws1.Range("G49").Value = '👀 À surveiller'
$ws1.Range("A50").Value = 'BANK OF AFRICA BF (BOABF)'
$ws1.Range("B50").Value = 0
$ws1.Range("C50").Value = 2
$ws1.Range("D50").Value = -4.37
$ws1.Range("E50").Value = -2.14
$ws1.Range("F50").Value = '🟡 Observer'
$ws1.Range("G50").Value = '➖ Neutre'
$ws1.Range("A51").Value = 'CFAO MOTORS CI (CFAC)'
$ws1.Range("B51").Value = 0
$ws1.Range("C51").Value = 2
$ws1.Range("D51").Value = -5.15
$ws1.Range("E51").Value = -2.21
$ws1.Range("F51").Value = '🟡 Observer'
$ws1.Range("G51").Value = '➖ Neutre'
$ws1.Range("A52").Value = 'BANK OF AFRICA BN (BOAB)'
$ws1.Range("B52").Value = 0
$ws1.Range("C52").Value = 2
$ws1.Range("D52").Value = -5.95
$ws1.Range("E52").Value = -2
$ws1.Range("F52").Value = '🟡 Observer'
$ws1.Range("G52").Value = '➖ Neutre'
$ws1.Range("A53").Value = 'BANK OF AFRICA NG (BOAN)'
$ws1.Range("B53").Value = 0
$ws1.Range("C53").Value = 2
$ws1.Range("D53").Value = -7.9
$ws1.Range("E53").Value = -1.25
$ws1.Range("F53").Value = '🟡 Observer'
$ws1.Range("G53").Value = '➖ Neutre'
$ws1.Range("A54").Value = 'SMB CI (SMBC)'
$ws1.Range("B54").Value = 0
$ws1.Range("C54").Value = 3
$ws1.Range("D54").Value = -9.42
$ws1.Range("E54").Value = -2.48
$ws1.Range("F54").Value = '🔴 Vente'
$ws1.Range("G54").Value = '⚠️ Risque de décrochage'

# --- Sheet2 "Top_YTD": update column B for rows 2-11 ---
$ws2.Range("B2").Value = 175295807.27
$ws2.Range("B3").Value = 14022040.04
$ws2.Range("B4").Value = 2997695.7
$ws2.Range("B5").Value = 2885302.52
$ws2.Range("B6").Value = 2560289.39
$ws2.Range("B7").Value = 1597505.19
$ws2.Range("B8").Value = 1507998.53
$ws2.Range("B9").Value = 1299956.91
$ws2.Range("B10").Value = 1056930.14
$ws2.Range("B11").Value = 240350.05

Write-Host "Update complete"
